$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46072
$ws.Range("B2").Value = 2.66
$ws.Range("C2").Value = 0.4
$ws.Range("D2").Value = 0.06
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.36
$ws.Range("H2").Value = 2.61
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7.46
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 0.52
$ws.Range("M2").Value = 0.99
$ws.Range("N2").Value = 0.99
$ws.Range("O2").Value = 0.99
$ws.Range("P2").Value = 0.99
$ws.Range("Q2").Value = 1.32
$ws.Range("R2").Value = 1.64
$ws.Range("S2").Value = 4.19
$ws.Range("T2").Value = 23.11
$ws.Range("U2").Value = 74.95
$ws.Range("V2").Value = 81.3
$ws.Range("W2").Value = 63.53
$ws.Range("X2").Value = 46.33
$ws.Range("Y2").Value = 33.4
$ws.Range("Z2").Value = 14.82
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 56.14
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 72.41
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 49.03
$ws.Range("AG2").Value = "0h-17h"
